# Add 2022-Q1 sheet (new fund holding data) and update the "总计" (Total)
# summary sheet with the new quarter's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Locate existing sheets we need as anchors / data sources.
# ---------------------------------------------------------------------------
$sheetQ4    = $wb.Worksheets.Item("2021-Q4")
$sheetTotal = $wb.Worksheets.Item("总计")

# We want the new "总计" sheet to end up with sheetId 5 (after the freshly
# created "2022-Q1" sheet takes sheetId 4), matching how the workbook was
# re-numbered upstream. Deleting and re-creating it (rather than just
# inserting a row) reproduces that numbering.
$sheetTotal.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Create the "2022-Q1" worksheet right after "2021-Q4".
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $sheetQ4)
$q1.Name = "2022-Q1"
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Style: header row + column A use the bold / centered / bordered look that
# is already used elsewhere in this workbook (e.g. on the "2021-Q4" sheet).
$sheetQ4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$sheetQ4.Range("A2:A9").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns B, D, E, F, G contain values that look numeric (fund codes,
# percentages, etc.) but must be stored as text, exactly like the rest of
# the workbook. Pre-formatting the range as Text keeps Excel from silently
# converting these to numbers when the values are assigned below.
$q1.Range("B2:B9").NumberFormat = "@"
$q1.Range("D2:G9").NumberFormat = "@"

$rows = @(
    @(0, "870009", "广发资管平衡精选一年持有混合A", "11.34", "94.29", "9.39", "1.0648", 2),
    @(1, "005449", "华夏行业龙头混合",               "11.19", "83.96", "2.92", "0.3267", 8),
    @(2, "011911", "华夏消费优选混合型证券投资基金A", "7.18",  "82.18", "4.13", "0.2965", 5),
    @(3, "872019", "广发资管平衡精选一年持有混合C", "1.54",  "94.29", "9.39", "0.1446", 2),
    @(4, "519678", "银河消费驱动混合",               "1.06",  "75.49", "8.66", "0.0918", 1),
    @(5, "519625", "银河君盛灵活配置混合A",          "4.91",  "20.05", "0.79", "0.0388", 10),
    @(6, "519626", "银河君盛灵活配置混合C",          "2.33",  "20.05", "0.79", "0.0184", 10),
    @(7, "011912", "华夏消费优选混合型证券投资基金C", "0.44",  "82.18", "4.13", "0.0182", 5)
)

$r = 2
foreach ($row in $rows) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    $q1.Range("G$r").Value = $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Re-create the "总计" worksheet right after "2022-Q1", with the new quarter
# inserted at the top of the summary table.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$sheetQ4.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$sheetQ4.Range("A2:A5").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalRows = @(
    @(0, "2022-Q1", 8,  2),
    @(1, "2021-Q4", 14, 5.72),
    @(2, "2021-Q3", 10, 3.18),
    @(3, "2021-Q2", 23, 5.35)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
    $r = $r + 1
}

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("2021-Q2").Activate()
